# Apply updated dSF (column F) values for the listed rows.
# Mirrors the "repull data, push all data, mean calculation" commit,
# which only changed the dSF column values for a subset of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    8  = -3
    9  = 0
    12 = -2
    13 = 3
    15 = -3
    19 = -3
    24 = -13
    26 = -4
    28 = -8
    30 = -8
    31 = -6
    32 = -9
    33 = 4
    38 = -2
    41 = -11
    44 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
